$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"   = -7.001199999999995
    "D14"  = -7.920899999999996
    "D16"  = -8.789900000000003
    "D21"  = -8.710999999999997
    "D23"  = -7.292199999999997
    "D25"  = -8.355599999999997
    "D26"  = -8.613999999999999
    "D29"  = -7.173499999999998
    "D40"  = -8.370999999999993
    "D53"  = -6.493899999999996
    "D57"  = -8.3721
    "D59"  = -8.203199999999997
    "D65"  = -7.942399999999996
    "D69"  = -7.139899999999998
    "D79"  = -6.144499999999999
    "D83"  = -8.9909
    "D91"  = -6.467999999999997
    "D93"  = -6.778699999999997
    "D100" = -8.231200000000003
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
